$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Cypher query text for the StatQuery row (shared by C2:C4)
$newQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Miniature Pinscher']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value2 = $newQuery
$ws.Range("C3").Value2 = $newQuery
$ws.Range("C4").Value2 = $newQuery

# Row heights shrink now that the query text is shorter
$ws.Rows.Item(2).RowHeight = 230.4
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# Reset the view: scroll back to top-left, zoom to 100%, keep B4 selected
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
$ws.Range("B4").Select()
